# Update column F values on the "展览" and "全部类型" sheets to reflect
# the latest generated data (commit: "Update gh-pages to output generated at 456a3b4").
#
# Changes (row -> old -> new):
#   3  : 3080 -> 3083
#   12 : 1394 -> 1395
#   16 : 52   -> 53
#   23 : 3262 -> 3263
#   25 : 160  -> 159
#   29 : 134  -> 135

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3  = 3083
    12 = 1395
    16 = 53
    23 = 3263
    25 = 159
    29 = 135
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
